# [IMP] Adjust CD Receivable Planning
# Rename the "Customer" / "Bank" filter labels to "Customer CD" / "Customer (bank)"
# on both sheets of the CD Receivable Planning report template.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("A7").Value = "Customer CD"
    $ws.Range("A8").Value = "Customer (bank)"
}
